$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell values are stored as text (coinranking.com price/volume feed uses
# "." as a thousands separator and keeps leading-zero/percent strings as text),
# so any new value that LOOKS like a plain number needs the cell pre-formatted
# as Text -- otherwise Excel.Range.Value silently re-types it as a Double.

$ws.Range("D2").Value = '40.019.84'
$ws.Range("E2").Value = '  +0.31%  '
$ws.Range("D3").Value = '2.212.95'
$ws.Range("E3").Value = '  -0.84%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '290.43'
$ws.Range("E5").Value = '  -3.04%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '88.45'
$ws.Range("E6").Value = '  +4.15%  '
$ws.Range("E7").Value = '  +0.15%  '
$ws.Range("E8").Value = '  -0.08%  '
$ws.Range("E9").Value = '  +0.61%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '30.79'
$ws.Range("E10").Value = '  +3.12%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0783'
$ws.Range("E11").Value = '  -0.07%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '47.69'
$ws.Range("E12").Value = '  +2.06%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.109'
$ws.Range("E13").Value = '  +1.99%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.46'
$ws.Range("E14").Value = '  +2.03%  '
$ws.Range("D15").Value = '2.556.41'
$ws.Range("E15").Value = '  -0.77%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.01'
$ws.Range("E16").Value = '  -1.19%  '
$ws.Range("D17").Value = '2.219.76'
$ws.Range("E17").Value = '  -0.41%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.728'
$ws.Range("E18").Value = '  +0.91%  '
$ws.Range("D19").Value = '39.949.05'
$ws.Range("E19").Value = '  +0.37%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.83'
$ws.Range("E20").Value = '  +12.82%  '
$ws.Range("D21").Value = '0.0₃0885'
$ws.Range("E21").Value = '  +0.58%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.81'
$ws.Range("E22").Value = '  +0.39%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '65.77'
$ws.Range("E23").Value = '  +0.74%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '235.91'
$ws.Range("E24").Value = '  +0.68%  '
$ws.Range("E25").Value = '  -0.06%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.46'
$ws.Range("E26").Value = '  +1.24%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.84'
$ws.Range("E27").Value = '  +1.05%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '22.58'
$ws.Range("E28").Value = '  -1.16%  '
$ws.Range("E29").Value = '  +4.36%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '9.25'
$ws.Range("E30").Value = '  +0.25%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '153.18'
$ws.Range("E31").Value = '  +2.20%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '32.20'
$ws.Range("E32").Value = '  -1.17%  '
$ws.Range("E33").Value = '  -0.13%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.96'
$ws.Range("E34").Value = '  +2.20%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0719'
$ws.Range("E35").Value = '  +2.27%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.39'
$ws.Range("E36").Value = '  -0.70%  '
$ws.Range("E37").Value = '  +6.34%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '16.06'
$ws.Range("E38").Value = '  -2.79%  '
$ws.Range("E39").Value = '  +0.45%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0999'
$ws.Range("E40").Value = '  +1.95%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.70'
$ws.Range("E41").Value = '  +2.21%  '
$ws.Range("D42").Value = '2.101.02'
$ws.Range("E42").Value = '  +8.59%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.82'
$ws.Range("E43").Value = '  +3.80%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.20'
$ws.Range("E44").Value = '  +2.36%  '
$ws.Range("E45").Value = '  +1.10%  '
$ws.Range("B46").Value = 'FraxShare'
$ws.Range("C46").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '9.90'
$ws.Range("E46").Value = '  +7.02%  '
$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '17.66'
$ws.Range("E47").Value = '  +6.93%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.66'
$ws.Range("E48").Value = '  +1.57%  '
$ws.Range("D49").Value = '2.430.65'
$ws.Range("E49").Value = '  -0.63%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '69.68'
$ws.Range("E50").Value = '  -2.49%  '
$ws.Range("B51").Value = 'Aave'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '88.74'
$ws.Range("E51").Value = '  -0.24%  '
